$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text would otherwise be auto-parsed as a number
# (e.g. "237.08", "1.000") are forced to Text format first so the literal
# string is preserved, matching the source inline-string cells.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.139.93"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "1.749.94"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "237.08"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").Value = "0.5542"
$ws.Range("E6").Value = "  +6.54%  "
$ws.Range("D7").Value = "0.9997"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "0.2846"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").Value = "0.06184"
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").Value = "1.748.51"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").Value = "0.07211"
$ws.Range("E11").Value = "  +2.68%  "
$ws.Range("D12").Value = "15.53"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "0.6518"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "4.656"
$ws.Range("E14").Value = "  +2.56%  "
$ws.Range("D15").Value = "78.45"
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").Value = "0.9991"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "26.023.10"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").Value = "11.79"
$ws.Range("E19").Value = "  +2.39%  "
$ws.Range("D20").Value = "0.000006784"
$ws.Range("E20").Value = "  +2.16%  "
$ws.Range("D21").Value = "1.971.29"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("D22").Value = "4.355"
$ws.Range("E22").Value = "  +4.57%  "
$ws.Range("D23").Value = "8.764"
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("D24").Value = "5.268"
$ws.Range("E24").Value = "  +2.17%  "
$ws.Range("D25").Value = "139.57"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").Value = "1.524"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("D27").Value = "15.35"
$ws.Range("E27").Value = "  +1.58%  "
$ws.Range("D28").Value = "1.812"
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("D29").Value = "105.70"
$ws.Range("E29").Value = "  +2.75%  "
$ws.Range("D30").Value = "0.08437"
$ws.Range("E30").Value = "  +1.66%  "
$ws.Range("D31").Value = "3.809"
$ws.Range("E31").Value = "  +3.71%  "
$ws.Range("D32").Value = "3.649"
$ws.Range("E32").Value = "  +5.82%  "
$ws.Range("D33").Value = "0.04645"
$ws.Range("E33").Value = "  +3.80%  "
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("D35").Value = "1.010"
$ws.Range("E35").Value = "  +2.09%  "
$ws.Range("D36").Value = "0.6315"
$ws.Range("E36").Value = "  +2.84%  "
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("D38").Value = "0.01620"
$ws.Range("E38").Value = "  +1.81%  "
$ws.Range("D39").Value = "1.981"
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("D40").Value = "0.9989"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").Value = "102.30"
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("D42").Value = "0.3949"
$ws.Range("E42").Value = "  +1.88%  "
$ws.Range("D43").Value = "0.7488"
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("D44").Value = "5.105"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").Value = "0.1154"
$ws.Range("E45").Value = "  +2.40%  "
$ws.Range("D46").Value = "6.374"
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("D47").Value = "0.05347"
$ws.Range("E47").Value = "  -2.24%  "
$ws.Range("D48").Value = "54.77"
$ws.Range("E48").Value = "  +3.26%  "
$ws.Range("D49").Value = "31.03"
$ws.Range("E49").Value = "  +3.27%  "
$ws.Range("D50").Value = "0.3504"
$ws.Range("E50").Value = "  +1.89%  "
$ws.Range("D51").Value = "7.622"
$ws.Range("E51").Value = "  -0.62%  "
